{"js": "// Week 3 status update:\n//   1. Fill in the GitHub URL that used to read \" N/A yet\".\n//   2. Collapse the split text runs in the Week 3 block back into single\n//      runs per paragraph (same visible text, consolidated formatting).\n\nconst body = context.document.body;\n\n// --- 1. GitHub URL -------------------------------------------------------\nconst urlResults = body.search(\"N/A yet\", { matchCase: true, matchWholeWord: false });\nurlResults.load(\"items\");\nawait context.sync();\n\nif (urlResults.items.length > 0) {\n  urlResults.items[0].insertText(\n    \"https://github.com/Jordan-m-jarvis/BoggleSolver\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- 2. Merge split runs in the Week 3 weekly-report block ---------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// These are the exact (already-correct) concatenated paragraph texts for\n// the Week 3 block; the source XML currently splits each of them across\n// two (or three) runs. Re-inserting the same text over each paragraph's\n// full range collapses it back down to one run, with no visible change.\nconst targetTexts = [\n  \"Number of hours worked this week: 12\",\n  \"Total number of hours worked on the project thus far:30\",\n  \"Number of total hours anticipated at completion:130\",\n  \"Accomplishments: Created colab notebook on google to power the object recognition model. Created a data-set of boggle boards with pictures and annotations on them. No GitHub yet, just experiments on the data-sets and attempting to get custom object recognition libraries running.\",\n  \"Challenges: Object recognition training takes around 25-40gb of ram when running. My laptop only has 16gb. I needed to tap into google colab to get enough ram and processing power to train the network. The problem is that sessions are only 1 hour long then you lose progress and have to restart. So I ruin the model for 55 min then download the partially trained model, then upload and continue training for another 55 min.\",\n  \"Plans / Goals for next week: Get gitHub up, put annotations and training data on github along with the trained or partially trained model and the code used to train the model. Get warp-affine working and a few filtering steps on the images taken from the webcam. Apply a grid and crop feature which takes the largest polygon and crops to it.\",\n  \"SPED Talk Insight (Briefly describe an insight or something interesting you learned from the SPED talks this week): I always struggled with grid layout not looking how I liked it. The SPED talk actually helped me understand why. Whenever I would make the grid I never added it to a smaller grid. Effectively I never figured out how to get it to automatically arrange in a way that I wanted it to\",\n];\nconst targetSet = new Set(targetTexts);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (targetSet.has(paragraph.text)) {\n    const range = paragraph.getRange();\n    range.insertText(paragraph.text, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Week 3 status update:\n#   1. Fill in the GitHub URL that used to read \" N/A yet\".\n#   2. Collapse the split text runs in the Week 3 block back into single\n#      runs per paragraph (same visible text, consolidated formatting).\n\n$d = $word.ActiveDocument\n\n# --- 1. GitHub URL --------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"N/A yet\"\n$find.Replacement.Text = \"https://github.com/Jordan-m-jarvis/BoggleSolver\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# --- 2. Merge split runs in the Week 3 weekly-report block ----------------\n# These are the exact (already-correct) concatenated paragraph texts for\n# the Week 3 block; the source XML currently splits each of them across\n# two (or three) runs. Deleting and re-inserting the same text over each\n# paragraph's range collapses it back down to one run, with no visible\n# change.\n$targets = @(\n  \"Number of hours worked this week: 12\",\n  \"Total number of hours worked on the project thus far:30\",\n  \"Number of total hours anticipated at completion:130\",\n  \"Accomplishments: Created colab notebook on google to power the object recognition model. Created a data-set of boggle boards with pictures and annotations on them. No GitHub yet, just experiments on the data-sets and attempting to get custom object recognition libraries running.\",\n  \"Challenges: Object recognition training takes around 25-40gb of ram when running. My laptop only has 16gb. I needed to tap into google colab to get enough ram and processing power to train the network. The problem is that sessions are only 1 hour long then you lose progress and have to restart. So I ruin the model for 55 min then download the partially trained model, then upload and continue training for another 55 min.\",\n  \"Plans / Goals for next week: Get gitHub up, put annotations and training data on github along with the trained or partially trained model and the code used to train the model. Get warp-affine working and a few filtering steps on the images taken from the webcam. Apply a grid and crop feature which takes the largest polygon and crops to it.\",\n  \"SPED Talk Insight (Briefly describe an insight or something interesting you learned from the SPED talks this week): I always struggled with grid layout not looking how I liked it. The SPED talk actually helped me understand why. Whenever I would make the grid I never added it to a smaller grid. Effectively I never figured out how to get it to automatically arrange in a way that I wanted it to\"\n)\n\nforeach ($p in $d.Paragraphs) {\n  $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($targets -contains $ptext) {\n    $r = $p.Range\n    $body = $d.Range($r.Start, $r.End - 1)\n    $body.Delete()\n    $body.InsertAfter($ptext)\n  }\n}\n"}
